# minor improvements to two PowerPoint slides
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18 ("Passes") - body placeholder: clarify that it is the definition
# of a "compiler pass" (not just "compiler") that some authors restrict to a
# traversal that involves disk I/O.
# ---------------------------------------------------------------------------
$slide18 = $p.Slides.Item(18)
$body18  = $slide18.Shapes.Item(4)
$tr18    = $body18.TextFrame.TextRange

$oldNote = "Note: Some authors restrict the definition of compiler to a traversal that involves disk I/O, but we will use a more general definition."
$newNote = "Note: Some authors restrict the definition of compiler pass to a traversal that involves disk I/O, but we will use a more general definition."

$fullText18 = $tr18.Text
$noteStart  = $fullText18.IndexOf($oldNote)
if ($noteStart -ge 0) {
    $noteRange = $tr18.Characters($noteStart + 1, $oldNote.Length)
    $noteRange.Text = $newNote
}

# ---------------------------------------------------------------------------
# Slide 19 ("Single-pass Versus Multi-pass Compilers") - body placeholder:
# reword the "ideal for multiprocessor systems" bullet.
# ---------------------------------------------------------------------------
$slide19 = $p.Slides.Item(19)
$body19  = $slide19.Shapes.Item(4)
$tr19    = $body19.TextFrame.TextRange

$oldBullet = "ideal for multiprocessor systems"
$newBullet = "can exploit concurrency and multiprocessor architectures"

$fullText19  = $tr19.Text
$bulletStart = $fullText19.IndexOf($oldBullet)
if ($bulletStart -ge 0) {
    $bulletRange = $tr19.Characters($bulletStart + 1, $oldBullet.Length)
    $bulletRange.Text = $newBullet
}
